$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -11.64929999999999
$ws.Range("B7").Value = 4.965999999999997
$ws.Range("A10").Value = -21.8295
$ws.Range("D10").Value = -7.9452
$ws.Range("A12").Value = -21.55210000000001
$ws.Range("D14").Value = -7.741700000000003
$ws.Range("B15").Value = 4.561799999999995
$ws.Range("A18").Value = -22.02560000000001
$ws.Range("C18").Value = -11.6268
$ws.Range("C19").Value = -11.48170000000001
$ws.Range("B20").Value = 8.5824
$ws.Range("C27").Value = -12.89539999999999
$ws.Range("B29").Value = 4.885100000000004
$ws.Range("B30").Value = 4.991700000000001
$ws.Range("B31").Value = 4.900300000000002
$ws.Range("D32").Value = -9.139799999999994
$ws.Range("D35").Value = -7.705300000000001
$ws.Range("A37").Value = -19.52919999999999
$ws.Range("B40").Value = 9.633799999999985
$ws.Range("C42").Value = -11.83410000000001
$ws.Range("D43").Value = -8.247499999999997
$ws.Range("C44").Value = -12.88119999999998
$ws.Range("C47").Value = -12.2209
$ws.Range("D49").Value = -8.024500000000003
$ws.Range("A55").Value = -22.0179
$ws.Range("D56").Value = -7.989999999999994
$ws.Range("C58").Value = -12.679
$ws.Range("A68").Value = -21.48800000000001
$ws.Range("B68").Value = 4.587099999999999
$ws.Range("D69").Value = -7.033499999999997
$ws.Range("C73").Value = -12.511
$ws.Range("B76").Value = 5.693899999999998
$ws.Range("A77").Value = -19.88849999999999
$ws.Range("A78").Value = -19.86239999999998
$ws.Range("D81").Value = -7.669099999999998
$ws.Range("B87").Value = 5.279999999999994
$ws.Range("B88").Value = 4.621299999999996
$ws.Range("D92").Value = -6.253899999999999
$ws.Range("C95").Value = -12.1236
$ws.Range("B96").Value = 5.381600000000005
$ws.Range("B98").Value = 5.447800000000004
$ws.Range("B101").Value = 9.758899999999993
$ws.Range("C101").Value = -12.49710000000001
$ws.Range("B102").Value = 8.540500000000003
